# Auto update stock data
# Update the "Date_1" column (A) for every data row from 2025/12/23 to
# 2025/12/24, and refresh a handful of EBITDA (column B) readings that
# changed alongside the date roll.
#
# NumberFormat is forced to Text ("@") right before each write so that the
# date-looking / number-looking strings are stored verbatim as text
# (matching the existing text cells) instead of being auto-converted into
# a date serial number or a float by Excel's normal type inference.
# ClearFormats() afterwards drops that scratch "@" formatting again so the
# cell's style/appearance is left exactly as it was before the edit - only
# the underlying text value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "2025/12/24"

$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    $cell = $ws.Range("A$r")
    $cell.NumberFormat = "@"
    $cell.Value = $newDate
    $cell.ClearFormats()
}

$ebitdaUpdates = @{
    8  = "8.54"
    14 = "3.01"
    26 = "11.24"
    32 = "27.95"
    44 = "11.24"
    50 = "11.55"
    56 = "31.76"
    62 = "11.71"
    68 = "13.25"
    74 = "16.66"
}

foreach ($r in $ebitdaUpdates.Keys) {
    $cell = $ws.Range("B$r")
    $cell.NumberFormat = "@"
    $cell.Value = $ebitdaUpdates[$r]
    $cell.ClearFormats()
}
